$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 131.27
$ws.Range("I15").Value = 131.27
$ws.Range("K15").Value = 393.8100000000001
$ws.Range("M15").Value = -224.8100000000001

$ws.Range("H64").Value = 4127.0386
$ws.Range("I64").Value = 3724.8572
$ws.Range("J64").Value = 4596.25
$ws.Range("K64").Value = 3724.8572
$ws.Range("L64").Value = 4596.25
$ws.Range("M64").Value = -3476.8572
$ws.Range("N64").Value = -5092.25

$ws.Range("H67").Value = 4127.0386
$ws.Range("I67").Value = 3724.8572
$ws.Range("J67").Value = 4596.25
$ws.Range("K67").Value = 3724.8572
$ws.Range("L67").Value = 4596.25
$ws.Range("M67").Value = -2866.8572
$ws.Range("N67").Value = -6312.25

$ws.Range("H74").Value = 5282.8237
$ws.Range("I74").Value = 4604.6665
$ws.Range("J74").Value = 5652.727
$ws.Range("K74").Value = 4604.6665
$ws.Range("L74").Value = 5652.727
$ws.Range("M74").Value = -3668.6665
$ws.Range("N74").Value = -7524.727

$ws.Range("H76").Value = 7548.567
$ws.Range("I76").Value = 8521.421
$ws.Range("J76").Value = 5868.1816
$ws.Range("K76").Value = 8521.421
$ws.Range("L76").Value = 5868.1816
$ws.Range("M76").Value = -8206.421
$ws.Range("N76").Value = -6498.1816

$ws.Range("H77").Value = 5282.8237
$ws.Range("I77").Value = 4604.6665
$ws.Range("J77").Value = 5652.727
$ws.Range("K77").Value = 23023.3325
$ws.Range("L77").Value = 28263.635
$ws.Range("M77").Value = -18343.3325
$ws.Range("N77").Value = -37623.63499999999

$ws.Range("H79").Value = 7548.567
$ws.Range("I79").Value = 8521.421
$ws.Range("J79").Value = 5868.1816
$ws.Range("K79").Value = 8521.421
$ws.Range("L79").Value = 5868.1816
$ws.Range("M79").Value = -7429.421
$ws.Range("N79").Value = -8052.1816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6153.227
$ws.Range("I63").Value = 7005.0713
$ws.Range("J63").Value = 4662.5
$ws.Range("K63").Value = 7005.0713
$ws.Range("L63").Value = 4662.5
$ws.Range("M63").Value = -6319.0713
$ws.Range("N63").Value = -6034.5

$ws.Range("H66").Value = 6153.227
$ws.Range("I66").Value = 7005.0713
$ws.Range("J66").Value = 4662.5
$ws.Range("K66").Value = 35025.35649999999
$ws.Range("L66").Value = 23312.5
$ws.Range("M66").Value = -31593.35649999999
$ws.Range("N66").Value = -30176.5

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19540.357
$ws.Range("I82").Value = 8040.625
$ws.Range("J82").Value = 34873.332
$ws.Range("K82").Value = 8040.625
$ws.Range("L82").Value = 34873.332
$ws.Range("M82").Value = -7657.625
$ws.Range("N82").Value = -35639.332

$ws.Range("H85").Value = 19540.357
$ws.Range("I85").Value = 8040.625
$ws.Range("J85").Value = 34873.332
$ws.Range("K85").Value = 8040.625
$ws.Range("L85").Value = 34873.332
$ws.Range("M85").Value = -6714.625
$ws.Range("N85").Value = -37525.332

$ws.Range("H86").Value = 7468.1177
$ws.Range("I86").Value = 5458.3076
$ws.Range("J86").Value = 14000
$ws.Range("K86").Value = 5458.3076
$ws.Range("L86").Value = 14000
$ws.Range("M86").Value = -4335.3076
$ws.Range("N86").Value = -16246

$ws.Range("H89").Value = 7468.1177
$ws.Range("I89").Value = 5458.3076
$ws.Range("J89").Value = 14000
$ws.Range("K89").Value = 27291.538
$ws.Range("L89").Value = 70000
$ws.Range("M89").Value = -21675.538
$ws.Range("N89").Value = -81232

$ws.Range("H94").Value = 1113.2572
$ws.Range("I94").Value = 1036.3043
$ws.Range("J94").Value = 1260.75
$ws.Range("K94").Value = 1036.3043
$ws.Range("L94").Value = 1260.75
$ws.Range("M94").Value = -585.3043
$ws.Range("N94").Value = -2162.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 29989.6
$ws.Range("J32").Value = 29989.6
$ws.Range("L32").Value = 29989.6
$ws.Range("N32").Value = -30581.6

$ws.Range("H70").Value = 9353.792
$ws.Range("I70").Value = 4312.8667
$ws.Range("J70").Value = 17755.334
$ws.Range("K70").Value = 4312.8667
$ws.Range("L70").Value = 17755.334
$ws.Range("M70").Value = -4042.8667
$ws.Range("N70").Value = -18295.334

$ws.Range("H73").Value = 9353.792
$ws.Range("I73").Value = 4312.8667
$ws.Range("J73").Value = 17755.334
$ws.Range("K73").Value = 4312.8667
$ws.Range("L73").Value = 17755.334
$ws.Range("M73").Value = -3376.8667
$ws.Range("N73").Value = -19627.334

$ws.Range("H80").Value = 4323.184
$ws.Range("I80").Value = 4972.926
$ws.Range("J80").Value = 2728.3635
$ws.Range("K80").Value = 4972.926
$ws.Range("L80").Value = 2728.3635
$ws.Range("M80").Value = -3974.926
$ws.Range("N80").Value = -4724.363499999999

$ws.Range("H83").Value = 4323.184
$ws.Range("I83").Value = 4972.926
$ws.Range("J83").Value = 2728.3635
$ws.Range("K83").Value = 24864.63
$ws.Range("L83").Value = 13641.8175
$ws.Range("M83").Value = -19872.63
$ws.Range("N83").Value = -23625.8175

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 9333
$ws.Range("J6").Value = 9333
$ws.Range("L6").Value = 9333
$ws.Range("N6").Value = -9557

$ws.Range("H74").Value = 32125
$ws.Range("J74").Value = 32125
$ws.Range("L74").Value = 32125
$ws.Range("N74").Value = -34121

$ws.Range("H77").Value = 32125
$ws.Range("J77").Value = 32125
$ws.Range("L77").Value = 96375
$ws.Range("N77").Value = -106359

$ws.Range("H82").Value = 2105.6743
$ws.Range("I82").Value = 2033.3928
$ws.Range("K82").Value = 2033.3928
$ws.Range("M82").Value = -1672.3928

$ws.Range("H85").Value = 2105.6743
$ws.Range("I85").Value = 2033.3928
$ws.Range("K85").Value = 2033.3928
$ws.Range("M85").Value = -785.3928000000001

$ws.Range("H93").Value = 693
$ws.Range("I93").Value = 726.8823
$ws.Range("J93").Value = 640.63635
$ws.Range("K93").Value = 726.8823
$ws.Range("L93").Value = 640.63635
$ws.Range("M93").Value = 521.1177
$ws.Range("N93").Value = -3136.63635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 24356.5
$ws.Range("J27").Value = 24356.5
$ws.Range("L27").Value = 24356.5
$ws.Range("N27").Value = -24494.5
